$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 5271.6665
$ws.Range("I12").Value = 5568.125
$ws.Range("K12").Value = 5568.125
$ws.Range("M12").Value = -5398.125
$ws.Range("H27").Value = 500
$ws.Range("J27").Value = 500
$ws.Range("L27").Value = 1500
$ws.Range("N27").Value = -1702
$ws.Range("H70").Value = 1758.8572
$ws.Range("I70").Value = 1499
$ws.Range("J70").Value = 1802.1666
$ws.Range("K70").Value = 4497
$ws.Range("L70").Value = 5406.4998
$ws.Range("M70").Value = -4227
$ws.Range("N70").Value = -5946.4998
$ws.Range("H73").Value = 1758.8572
$ws.Range("I73").Value = 1499
$ws.Range("J73").Value = 1802.1666
$ws.Range("K73").Value = 4497
$ws.Range("L73").Value = 5406.4998
$ws.Range("M73").Value = -3561
$ws.Range("N73").Value = -7278.4998
$ws.Range("H80").Value = 306.22223
$ws.Range("I80").Value = 195.28572
$ws.Range("J80").Value = 376.81818
$ws.Range("K80").Value = 585.85716
$ws.Range("L80").Value = 1130.45454
$ws.Range("M80").Value = 412.14284
$ws.Range("N80").Value = -3126.45454
$ws.Range("H83").Value = 306.22223
$ws.Range("I83").Value = 195.28572
$ws.Range("J83").Value = 376.81818
$ws.Range("K83").Value = 1757.57148
$ws.Range("L83").Value = 3391.36362
$ws.Range("M83").Value = 3234.42852
$ws.Range("N83").Value = -13375.36362
$ws.Range("H86").Value = 1895.2609
$ws.Range("I86").Value = 2125.2856
$ws.Range("J86").Value = 1537.4445
$ws.Range("K86").Value = 2125.2856
$ws.Range("L86").Value = 1537.4445
$ws.Range("M86").Value = -1002.2856
$ws.Range("N86").Value = -3783.4445
$ws.Range("H89").Value = 1895.2609
$ws.Range("I89").Value = 2125.2856
$ws.Range("J89").Value = 1537.4445
$ws.Range("K89").Value = 10626.428
$ws.Range("L89").Value = 7687.2225
$ws.Range("M89").Value = -5010.428
$ws.Range("N89").Value = -18919.2225
$ws.Range("H106").Value = 12610.6
$ws.Range("I106").Value = 2941.75
$ws.Range("K106").Value = 2941.75
$ws.Range("M106").Value = -2310.75
$ws.Range("H113").Value = 4776.4707
$ws.Range("I113").Value = 4457.5713
$ws.Range("K113").Value = 4457.5713
$ws.Range("M113").Value = -1203.5713
$ws.Range("H116").Value = 4194.0835
$ws.Range("I116").Value = 3991.2
$ws.Range("K116").Value = 3991.2
$ws.Range("M116").Value = -549.1999999999998
$ws.Range("H138").Value = 2252.1187
$ws.Range("I138").Value = 963.40625
$ws.Range("J138").Value = 3779.4814
$ws.Range("K138").Value = 2890.21875
$ws.Range("L138").Value = 11338.4442
$ws.Range("M138").Value = 2249.78125
$ws.Range("N138").Value = -21618.4442
$ws.Range("H141").Value = 26379.432
$ws.Range("I141").Value = 26379.432
$ws.Range("K141").Value = 79138.296
$ws.Range("M141").Value = -73958.296

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 250875.25
$ws.Range("I4").Value = 250875.25
$ws.Range("K4").Value = 250875.25
$ws.Range("M4").Value = -250759.25
$ws.Range("H32").Value = 40981.332
$ws.Range("I32").Value = 24717.906
$ws.Range("K32").Value = 24717.906
$ws.Range("M32").Value = -24430.906
$ws.Range("H37").Value = 8347472.5
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H45").Value = 462058.8
$ws.Range("I45").Value = 563703.25
$ws.Range("K45").Value = 563703.25
$ws.Range("M45").Value = -563326.25
$ws.Range("H55").Value = 19989.75
$ws.Range("H61").Value = 1896
$ws.Range("I61").Value = 1929.1428
$ws.Range("K61").Value = 1929.1428
$ws.Range("M61").Value = -1717.1428
$ws.Range("H74").Value = 2005.1666
$ws.Range("I74").Value = 2011.7693
$ws.Range("K74").Value = 2011.7693
$ws.Range("M74").Value = -1137.7693
$ws.Range("H77").Value = 2005.1666
$ws.Range("I77").Value = 2011.7693
$ws.Range("K77").Value = 10058.8465
$ws.Range("M77").Value = -5690.8465
$ws.Range("H80").Value = 19998.334
$ws.Range("J80").Value = 19998.334
$ws.Range("L80").Value = 19998.334
$ws.Range("N80").Value = -21994.334
$ws.Range("H83").Value = 19998.334
$ws.Range("J83").Value = 19998.334
$ws.Range("L83").Value = 59995.00199999999
$ws.Range("N83").Value = -69979.00199999999
$ws.Range("H136").Value = 1896
$ws.Range("I136").Value = 1929.1428
$ws.Range("K136").Value = 5787.428400000001
$ws.Range("M136").Value = -3237.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5943.357
$ws.Range("I105").Value = 7318.8184
$ws.Range("K105").Value = 7318.8184
$ws.Range("M105").Value = -5571.8184
$ws.Range("H134").Value = 1042.3158
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2390.2285
$ws.Range("I134").Value = 2530.1155
$ws.Range("J134").Value = 1986.1111
$ws.Range("K134").Value = 7590.3465
$ws.Range("L134").Value = 5958.3333
$ws.Range("M134").Value = -5055.3465
$ws.Range("N134").Value = -11028.3333
$ws.Range("H135").Value = 89993.336
$ws.Range("J135").Value = 89993.336
$ws.Range("L135").Value = 89993.336
$ws.Range("N135").Value = -100133.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2997
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2761.3667
$ws.Range("I132").Value = 1997.1666
$ws.Range("J132").Value = 2952.4167
$ws.Range("K132").Value = 17974.4994
$ws.Range("L132").Value = 26571.7503
$ws.Range("M132").Value = -15444.4994
$ws.Range("N132").Value = -31631.7503
$ws.Range("H134").Value = 3658.9
$ws.Range("I134").Value = 1823.625
$ws.Range("K134").Value = 5470.875
$ws.Range("M134").Value = -400.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4091.5454
$ws.Range("I102").Value = 4100.7
$ws.Range("K102").Value = 4100.7
$ws.Range("M102").Value = -2478.7
$ws.Range("H109").Value = 200285
$ws.Range("J109").Value = 200285
$ws.Range("L109").Value = 200285
$ws.Range("N109").Value = -202365
$ws.Range("H126").Value = 3841.8572
$ws.Range("I126").Value = 3740.2
$ws.Range("J126").Value = 4096
$ws.Range("K126").Value = 11220.6
$ws.Range("L126").Value = 12288
$ws.Range("M126").Value = -8750.599999999999
$ws.Range("N126").Value = -17228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18222.277
$ws.Range("I7").Value = 24001
$ws.Range("J7").Value = 6664.8335
$ws.Range("K7").Value = 24001
$ws.Range("L7").Value = 6664.8335
$ws.Range("M7").Value = -23889
$ws.Range("N7").Value = -6888.8335
$ws.Range("H40").Value = 6155.2666
$ws.Range("I40").Value = 6028.6665
$ws.Range("K40").Value = 6028.6665
$ws.Range("M40").Value = -5892.6665
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H126").Value = 18222.277
$ws.Range("I126").Value = 24001
$ws.Range("J126").Value = 6664.8335
$ws.Range("K126").Value = 72003
$ws.Range("L126").Value = 19994.5005
$ws.Range("M126").Value = -69533
$ws.Range("N126").Value = -24934.5005
$ws.Range("H132").Value = 4260.44
$ws.Range("I132").Value = 4098.4707
$ws.Range("J132").Value = 4604.625
$ws.Range("K132").Value = 12295.4121
$ws.Range("L132").Value = 13813.875
$ws.Range("M132").Value = -9765.4121
$ws.Range("N132").Value = -18873.875
$ws.Range("H136").Value = 4544.2144
$ws.Range("I136").Value = 3233.4285
$ws.Range("K136").Value = 9700.2855
$ws.Range("M136").Value = -7150.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 8096.05
$ws.Range("I126").Value = 4785.4287
$ws.Range("K126").Value = 14356.2861
$ws.Range("M126").Value = -11886.2861
$ws.Range("H132").Value = 10795.192
$ws.Range("I132").Value = 11331.25
$ws.Range("K132").Value = 33993.75
$ws.Range("M132").Value = -31463.75
$ws.Range("H136").Value = 1230.8334
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 1461.6666
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 4384.9998
$ws.Range("M136").Value = -450
